# poeple import and people重点人员识别记录
# - Rename sheet from "侦码后台信息" to "人口库导入"
# - Switch calculation mode from manual to automatic
# - Move the saved selection/active cell from E18 to C31

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet was renamed to reflect the people-import purpose of the workbook.
$ws.Name = "人口库导入"

# Calculation was left in manual mode before; switch back to automatic.
$excel.Calculation = -4105

# Update the active cell / selection recorded in the sheet view.
$ws.Range("C31").Select()
